# Update betting odds on Sheet1 to reflect the latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3
$ws.Range("N3").Value = 2.25
$ws.Range("O3").Value = 1.62

# Row 6
$ws.Range("G6").Value = 1.55
$ws.Range("H6").Value = 4.6
$ws.Range("I6").Value = 4.55
$ws.Range("M6").Value = 6.7
$ws.Range("Q6").Value = 4.55
$ws.Range("R6").Value = 1.34
$ws.Range("S6").Value = 3
$ws.Range("T6").Value = 17
$ws.Range("W6").Value = 16
$ws.Range("X6").Value = 11
$ws.Range("Y6").Value = 14
$ws.Range("Z6").Value = 37
$ws.Range("AA6").Value = 11.75
$ws.Range("AB6").Value = 12
$ws.Range("AC6").Value = 27
$ws.Range("AD6").Value = 100
$ws.Range("AG6").Value = 16.5
$ws.Range("AJ6").Value = 25

# Row 7
$ws.Range("I7").Value = 1.77
$ws.Range("J7").Value = 1.06
$ws.Range("K7").Value = 10
$ws.Range("N7").Value = 2.08
$ws.Range("O7").Value = 1.73

# Row 8
$ws.Range("G8").Value = 2.1
$ws.Range("I8").Value = 3.4
$ws.Range("R8").Value = 1.95
$ws.Range("S8").Value = 1.8
$ws.Range("U8").Value = 9.5
$ws.Range("V8").Value = 9.5
$ws.Range("AG8").Value = 12
$ws.Range("AI8").Value = 29
